$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1597
$ws.Range("F4").Value = 843
$ws.Range("F6").Value = 70
$ws.Range("F7").Value = 1149
$ws.Range("F8").Value = 752
$ws.Range("F9").Value = 799
$ws.Range("F10").Value = 1449
$ws.Range("F15").Value = 194
$ws.Range("F17").Value = 476
$ws.Range("F18").Value = 28
$ws.Range("F19").Value = 27
$ws.Range("F20").Value = 4
$ws.Range("F22").Value = 298
$ws.Range("F23").Value = 554
$ws.Range("F25").Value = 761
$ws.Range("F27").Value = 182
$ws.Range("F28").Value = 371

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 1002
$ws.Range("F5").Value = 270
$ws.Range("F10").Value = 81

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 239

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 239
$ws.Range("F4").Value = 1597
$ws.Range("F6").Value = 843
$ws.Range("F8").Value = 1002
$ws.Range("F9").Value = 70
$ws.Range("F10").Value = 1149
$ws.Range("F11").Value = 752
$ws.Range("F12").Value = 799
$ws.Range("F13").Value = 1449
$ws.Range("F18").Value = 194
$ws.Range("F20").Value = 476
$ws.Range("F21").Value = 28
$ws.Range("F22").Value = 27
$ws.Range("F24").Value = 4
$ws.Range("F25").Value = 270
$ws.Range("F27").Value = 298
$ws.Range("F31").Value = 554
$ws.Range("F32").Value = 570
$ws.Range("F33").Value = 761
$ws.Range("F36").Value = 182
$ws.Range("F38").Value = 81
$ws.Range("F39").Value = 81
$ws.Range("F41").Value = 371
